$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain text (matches original inlineStr text cells)
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '70.679.21'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').Value = '3.581.45'
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '597.68'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = '173.23'
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('D7').Value = '3.578.28'
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  +5.65%  '
$ws.Range('E11').Value = '  +7.69%  '
$ws.Range('D12').Value = '0.589'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').Value = '46.71'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').Value = '0.0000278'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').Value = '4.159.20'
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '614.33'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = '3.584.54'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('D19').Value = '70.751.00'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').Value = '17.47'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = '0.885'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  -16.75%  '
$ws.Range('D24').Value = '15.84'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').Value = '97.06'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '2.63'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').Value = '33.95'
$ws.Range('E29').Value = '  +4.00%  '
$ws.Range('D30').Value = '9.19'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = '8.38'
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('D32').Value = '3.06'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').Value = '7.21'
$ws.Range('E33').Value = '  +4.34%  '
$ws.Range('D34').Value = '645.88'
$ws.Range('E34').Value = '  +2.63%  '
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').Value = '3.68'
$ws.Range('E36').Value = '  +5.99%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').Value = '10.82'
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').Value = '0.0479'
$ws.Range('E39').Value = '  +7.49%  '
$ws.Range('D40').Value = '57.36'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +5.00%  '
$ws.Range('D43').Value = '3.393.44'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '0.323'
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').Value = '0.0₃0712'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').Value = '32.91'
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E47').Value = '  +6.18%  '
$ws.Range('D48').Value = '2.65'
$ws.Range('E48').Value = '  +4.61%  '
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').Value = '132.98'
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('E51').Value = '  -0.09%  '
